$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D2: '27.367.85' -> '27.308.10'; E2: '  +2.64%  ' -> '  +2.11%  '
$ws.Range("D2").Value = "27.308.10"
$ws.Range("E2").Value = "  +2.11%  "

# Row 3: D3: '1.822.20' -> '1.819.35'; E3: '  +1.66%  ' -> '  +1.30%  '
$ws.Range("D3").Value = "1.819.35"
$ws.Range("E3").Value = "  +1.30%  "

# Row 4: D4: '0.9989' -> '1.001'; E4: '  -0.38%  ' -> '  +0.03%  '
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.03%  "

# Row 5: D5: '313.85' -> '313.22'; E5: '  +1.76%  ' -> '  +1.45%  '
$ws.Range("D5").Value = "'313.22"
$ws.Range("E5").Value = "  +1.45%  "

# Row 6: D6: '0.9997' -> '1.000'; E6: '  -0.22%  ' -> '  -0.08%  '
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.08%  "

# Row 7: D7: '0.4659' -> '0.4649'; E7: '  +5.67%  ' -> '  +5.51%  '
$ws.Range("D7").Value = "'0.4649"
$ws.Range("E7").Value = "  +5.51%  "

# Row 8: D8: '0.3773' -> '0.3766'; E8: '  +2.98%  ' -> '  +2.58%  '
$ws.Range("D8").Value = "'0.3766"
$ws.Range("E8").Value = "  +2.58%  "

# Row 9: D9: '0.07441' -> '0.07421'; E9: '  +1.61%  ' -> '  +1.16%  '
$ws.Range("D9").Value = "'0.07421"
$ws.Range("E9").Value = "  +1.16%  "

# Row 10: D10: '0.8722' -> '0.8701'; E10: '  +2.09%  ' -> '  +1.68%  '
$ws.Range("D10").Value = "'0.8701"
$ws.Range("E10").Value = "  +1.68%  "

# Row 11: D11: '20.68' -> '20.62'; E11: '  +0.71%  ' -> '  +0.14%  '
$ws.Range("D11").Value = "'20.62"
$ws.Range("E11").Value = "  +0.14%  "

# Row 12: D12: '1.825.53' -> '1.817.89'; E12: '  +1.80%  ' -> '  +1.13%  '
$ws.Range("D12").Value = "1.817.89"
$ws.Range("E12").Value = "  +1.13%  "

# Row 13: D13: '6.685' -> '6.661'; E13: '  +1.35%  ' -> '  +0.70%  '
$ws.Range("D13").Value = "'6.661"
$ws.Range("E13").Value = "  +0.70%  "

# Row 14: D14: '5.414' -> '5.399'; E14: '  +3.27%  ' -> '  +2.64%  '
$ws.Range("D14").Value = "'5.399"
$ws.Range("E14").Value = "  +2.64%  "

# Row 15: D15: '0.07103' -> '0.07113'; E15: '  +0.70%  ' -> '  +0.73%  '
$ws.Range("D15").Value = "'0.07113"
$ws.Range("E15").Value = "  +0.73%  "

# Row 16: D16: '92.23' -> '92.08'; E16: '  +1.38%  ' -> '  +0.82%  '
$ws.Range("D16").Value = "'92.08"
$ws.Range("E16").Value = "  +0.82%  "

# Row 17: D17: '0.9998' -> '1.002'; E17: '  -0.36%  ' -> '  +0.00%  '
$ws.Range("D17").Value = "'1.002"
$ws.Range("E17").Value = "  +0.00%  "

# Row 18: D18: '0.000008776' -> '0.000008756'; E18: '  +1.40%  ' -> '  +1.19%  '
$ws.Range("D18").Value = "'0.000008756"
$ws.Range("E18").Value = "  +1.19%  "

# Row 19: E19: '  -0.12%  ' -> '  -0.19%  '
$ws.Range("E19").Value = "  -0.19%  "

# Row 20: D20: '14.96' -> '14.94'; E20: '  +1.52%  ' -> '  +1.35%  '
$ws.Range("D20").Value = "'14.94"
$ws.Range("E20").Value = "  +1.35%  "

# Row 21: D21: '27.371.22' -> '27.333.32'; E21: '  +2.56%  ' -> '  +2.25%  '
$ws.Range("D21").Value = "27.333.32"
$ws.Range("E21").Value = "  +2.25%  "

# Row 22: D22: '5.310' -> '5.302'; E22: '  +3.25%  ' -> '  +3.00%  '
$ws.Range("D22").Value = "'5.302"
$ws.Range("E22").Value = "  +3.00%  "

# Row 23: D23: '10.92' -> '10.90'; E23: '  +1.31%  ' -> '  +0.92%  '
$ws.Range("D23").Value = "'10.90"
$ws.Range("E23").Value = "  +0.92%  "

# Row 24: D24: '2.052.48' -> '2.053.23'; E24: '  +1.78%  ' -> '  +1.90%  '
$ws.Range("D24").Value = "2.053.23"
$ws.Range("E24").Value = "  +1.90%  "

# Row 25: D25: '1.940' -> '1.941'; E25: '  -1.66%  ' -> '  -1.70%  '
$ws.Range("D25").Value = "'1.941"
$ws.Range("E25").Value = "  -1.70%  "

# Row 26: D26: '151.56' -> '151.49'; E26: '  -0.03%  ' -> '  -0.07%  '
$ws.Range("D26").Value = "'151.49"
$ws.Range("E26").Value = "  -0.07%  "

# Row 27: D27: '2.261' -> '2.259'; E27: '  +4.16%  ' -> '  +3.30%  '
$ws.Range("D27").Value = "'2.259"
$ws.Range("E27").Value = "  +3.30%  "

# Row 28: D28: '18.63' -> '18.60'; E28: '  +1.54%  ' -> '  +1.25%  '
$ws.Range("D28").Value = "'18.60"
$ws.Range("E28").Value = "  +1.25%  "

# Row 29: D29: '5.302' -> '5.287'; E29: '  +3.11%  ' -> '  +2.44%  '
$ws.Range("D29").Value = "'5.287"
$ws.Range("E29").Value = "  +2.44%  "

# Row 30: D30: '117.12' -> '117.08'; E30: '  +0.25%  ' -> '  -0.07%  '
$ws.Range("D30").Value = "'117.08"
$ws.Range("E30").Value = "  -0.07%  "

# Row 31: D31: '0.08899' -> '0.08910'; E31: '  +1.68%  ' -> '  +1.48%  '
$ws.Range("D31").Value = "'0.08910"
$ws.Range("E31").Value = "  +1.48%  "

# Row 32: D32: '0.7837' -> '0.7797'; E32: '  +6.98%  ' -> '  +6.33%  '
$ws.Range("D32").Value = "'0.7797"
$ws.Range("E32").Value = "  +6.33%  "

# Row 33: D33: '1.185' -> '1.182'; E33: '  +3.89%  ' -> '  +3.23%  '
$ws.Range("D33").Value = "'1.182"
$ws.Range("E33").Value = "  +3.23%  "

# Row 34: D34: '4.536' -> '4.528'; E34: '  +2.74%  ' -> '  +2.34%  '
$ws.Range("D34").Value = "'4.528"
$ws.Range("E34").Value = "  +2.34%  "

# Row 35: D35: '2.925' -> '2.921'; E35: '  +0.71%  ' -> '  +0.72%  '
$ws.Range("D35").Value = "'2.921"
$ws.Range("E35").Value = "  +0.72%  "

# Row 36: D36: '0.9991' -> '0.9997'; E36: '  -0.31%  ' -> '  -0.09%  '
$ws.Range("D36").Value = "'0.9997"
$ws.Range("E36").Value = "  -0.09%  "

# Row 37: D37: '1.098' -> '1.102'; E37: '  +1.22%  ' -> '  +1.75%  '
$ws.Range("D37").Value = "'1.102"
$ws.Range("E37").Value = "  +1.75%  "

# Row 38: D38: '0.01974' -> '0.01970'; E38: '  +1.61%  ' -> '  +1.25%  '
$ws.Range("D38").Value = "'0.01970"
$ws.Range("E38").Value = "  +1.25%  "

# Row 39: D39: '0.05265' -> '0.05260'; E39: '  +2.14%  ' -> '  +1.90%  '
$ws.Range("D39").Value = "'0.05260"
$ws.Range("E39").Value = "  +1.90%  "

# Row 40: D40: '7.295' -> '7.280'; E40: '  +5.17%  ' -> '  +4.64%  '
$ws.Range("D40").Value = "'7.280"
$ws.Range("E40").Value = "  +4.64%  "

# Row 41: D41: '2.389' -> '2.382'; E41: '  +21.99%  ' -> '  +21.67%  '
$ws.Range("D41").Value = "'2.382"
$ws.Range("E41").Value = "  +21.67%  "

# Row 42: B42: 'TheSandbox' -> 'MXToken'; C42: 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand' -> 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D42: '0.5309' -> '2.912'; E42: '  +1.93%  ' -> '  +3.78%  '
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "'2.912"
$ws.Range("E42").Value = "  +3.78%  "

# Row 43: B43: 'MXToken' -> 'TheSandbox'; C43: 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' -> 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D43: '2.904' -> '0.5295'; E43: '  +3.41%  ' -> '  +1.53%  '
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.5295"
$ws.Range("E43").Value = "  +1.53%  "

# Row 44: D44: '0.1691' -> '0.1688'; E44: '  +1.15%  ' -> '  +0.70%  '
$ws.Range("D44").Value = "'0.1688"
$ws.Range("E44").Value = "  +0.70%  "

# Row 45: D45: '8.635' -> '8.606'; E45: '  +2.95%  ' -> '  +2.34%  '
$ws.Range("D45").Value = "'8.606"
$ws.Range("E45").Value = "  +2.34%  "

# Row 46: D46: '0.5058' -> '0.5046'; E46: '  +0.99%  ' -> '  +0.43%  '
$ws.Range("D46").Value = "'0.5046"
$ws.Range("E46").Value = "  +0.43%  "

# Row 47: D47: '10.53' -> '10.51'; E47: '  +1.52%  ' -> '  +1.58%  '
$ws.Range("D47").Value = "'10.51"
$ws.Range("E47").Value = "  +1.58%  "

# Row 48: D48: '105.69' -> '105.45'; E48: '  +0.95%  ' -> '  +0.59%  '
$ws.Range("D48").Value = "'105.45"
$ws.Range("E48").Value = "  +0.59%  "

# Row 49: D49: '1.678' -> '1.675'; E49: '  +1.43%  ' -> '  +1.07%  '
$ws.Range("D49").Value = "'1.675"
$ws.Range("E49").Value = "  +1.07%  "

# Row 50: D50: '0.9991' -> '0.9986'; E50: '  -0.18%  ' -> '  -0.26%  '
$ws.Range("D50").Value = "'0.9986"
$ws.Range("E50").Value = "  -0.26%  "

# Row 51: D51: '0.06339' -> '0.06334'; E51: '  +0.88%  ' -> '  +0.80%  '
$ws.Range("D51").Value = "'0.06334"
$ws.Range("E51").Value = "  +0.80%  "
